# Add an "affectedFuels" worksheet, carrying the "...among affected" fuel
# share parameter rows that used to live at the bottom of "main", out into
# their own tab (right after "main").
#
# feat: added affected fuel category to advanced parameters

$wb = $excel.ActiveWorkbook

$mainSheet = $wb.Worksheets.Item("main")

# Duplicate "main" (keeps column widths / styles / header row intact) and
# drop it immediately after "main" itself.
$mainSheet.Copy($null, $mainSheet) | Out-Null
$affectedFuels = $wb.Worksheets.Item($mainSheet.Index + 1)
$affectedFuels.Name = "affectedFuels"

# The new sheet only needs to keep the header row (row 1) plus the
# "Share of <fuel> among affected" rows (rows 6-12 in the original "main"
# sheet) - remove the rows in between (the old rows 2-5).
$affectedFuels.Rows("2:5").Delete() | Out-Null

# Back on "main", those "Share of <fuel> among affected" rows now live only
# on "affectedFuels", so remove them from "main".
$mainSheet.Rows("6:12").Select() | Out-Null
$mainSheet.Rows("6:12").Delete() | Out-Null

# Make the freshly added sheet the active tab, matching the saved view.
$affectedFuels.Activate() | Out-Null
$affectedFuels.Range("B14").Select() | Out-Null
